# Added predicted category for deviations and improvements points (ML)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F ("Forbedringspunkter" and everything
# to its right shifts one column to the right, from F:T -> G:U).
$ws.Columns("F").Insert()

# New column F holds the ML-predicted category for the deviation (avvik).
$ws.Cells.Item(1, 6).Value = "Kategori (avvik)"
$ws.Cells.Item(2, 6).Value = "['teknisksikkerhet']"
$ws.Cells.Item(3, 6).Value = "['styring']"
$ws.Cells.Item(4, 6).Value = "['styring']"
$ws.Cells.Item(5, 6).Value = "['materialhåndtering']"
$ws.Cells.Item(6, 6).Value = "['styring']"
$ws.Cells.Item(7, 6).Value = "['styring']"

# Column K (previously J, "Kategori (forbedringer)") now holds the
# ML-predicted category for the improvement point instead of the old
# category_prediction(test_imp.description) placeholder text.
$ws.Cells.Item(2, 11).Value = "['styring']"
$ws.Cells.Item(3, 11).Value = "['teknisksikkerhet']"
